$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '40.212.77'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +3.39%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.245.74'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +1.02%  '
$ws.Range("E4").Value = '  -0.03%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '295.07'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.47%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '87.06'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +9.30%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.518'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +2.36%  '
$ws.Range("E8").Value = '  -0.01%  '
$ws.Range("E9").Value = '  +4.23%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '31.33'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +12.60%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0801'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +3.80%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '47.26'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +2.38%  '
$ws.Range("E13").Value = '  +1.21%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '6.49'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +6.60%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '2.595.76'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +1.23%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '14.28'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +2.33%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '2.253.13'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +0.90%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.740'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +4.02%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '40.141.25'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +3.36%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.0₃0896'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +4.41%  '
$ws.Range("E21").Value = '  +2.50%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '10.69'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +8.89%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '65.84'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +1.80%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '236.73'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +5.36%  '
$ws.Range("E25").Value = '  +0.13%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.48'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +3.78%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '1.85'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +8.07%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '23.08'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +4.34%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.24'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +2.38%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '9.30'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +5.12%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '33.45'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +7.85%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '154.51'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +3.95%  '
$ws.Range("E33").Value = '  -0.12%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '4.92'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +3.18%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.0721'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +5.68%  '
$ws.Range("E36").Value = '  +3.22%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '16.63'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +16.08%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.101'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +6.34%  '
$ws.Range("E39").Value = '  +3.17%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.73'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +3.24%  '
$ws.Range("E41").Value = '  +7.90%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '3.81'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +6.45%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '2.024.99'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +6.65%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '2.22'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +11.33%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.0273'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +7.96%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '9.96'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +11.29%  '
$ws.Range("E47").Value = '  +2.46%  '
$ws.Range("E48").Value = '  +3.65%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '2.476.29'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +1.69%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '71.80'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +5.13%  '
$ws.Range("E51").Value = '  +16.61%  '
